# Add a new row (13) to the Maven guide worksheet documenting the
# "mvn package" shortcut scenario and its result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 values -------------------------------------------------------
$ws.Range("A13").Value = 12

$text1 = "We can also use command mvn package directly from the begining it will execute the clean, compile, test and package at a time and it will run one by one automatically. First will clear the environment and will create the target file and package it directly"
$ws.Range("B13").Value = $text1

$text2 = "Process done automatically and created a Target folder and showed Build success."
$ws.Range("C13").Value = $text2

# --- Rich text (bold) formatting for B13 ----------------------------------
# Runs:
#  1: "We can also use command "        (plain)
#  2: "mvn package directly "           (bold)
#  3: "from the begining"               (plain)
#  4: " it will"                        (bold)
#  5: " execute the "                   (plain)
#  6: "clean, compile, test and package " (bold)
#  7: "at a time and it will run one by one automatically. First will clear the environment and will create the target file and package it directly" (plain)

$r = $ws.Range("B13").Characters(25, 21)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("B13").Characters(46, 17)
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("B13").Characters(63, 8)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("B13").Characters(71, 13)
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("B13").Characters(84, 33)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("B13").Characters(117, 140)
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

# --- Rich text (bold) formatting for C13 ----------------------------------
# Runs:
#  1: "Process done automatically and created a " (plain)
#  2: "Target folder "                            (bold)
#  3: "and showed "                                (plain)
#  4: "Build success."                             (bold)

$r = $ws.Range("C13").Characters(42, 14)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("C13").Characters(56, 11)
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

$r = $ws.Range("C13").Characters(67, 14)
$r.Font.Bold = $true
$r.Font.Size = 11
$r.Font.Name = "Aptos Narrow"

# --- Row height (matches the autofit height Excel computed for wrapped text)
$ws.Rows.Item(13).RowHeight = 57.6

# --- Update the view: scroll & selection as left by the author ----------
$null = $ws.Range("C13").Select()

Write-Host "Row 13 added."
